{"js": "// Replace each two-digit-by-two-digit multiplication problem text in the\n// document body with its updated value, as described by the diff.\nconst replacements = [\n  [\"20\u00d765=1300\", \"21\u00d776=1596\"],\n  [\"99\u00d785=8415\", \"54\u00d791=4914\"],\n  [\"59\u00d790=5310\", \"38\u00d787=3306\"],\n  [\"15\u00d751=765\", \"61\u00d758=3538\"],\n  [\"89\u00d787=7743\", \"43\u00d714=602\"],\n  [\"96\u00d758=5568\", \"36\u00d748=1728\"],\n  [\"22\u00d791=2002\", \"25\u00d798=2450\"],\n  [\"42\u00d797=4074\", \"43\u00d785=3655\"],\n  [\"25\u00d749=1225\", \"24\u00d722=528\"],\n  [\"27\u00d774=1998\", \"25\u00d716=400\"],\n  [\"29\u00d759=1711\", \"74\u00d790=6660\"],\n  [\"94\u00d719=1786\", \"27\u00d731=837\"],\n  [\"21\u00d797=2037\", \"13\u00d753=689\"],\n  [\"65\u00d735=2275\", \"91\u00d731=2821\"],\n  [\"63\u00d795=5985\", \"65\u00d775=4875\"],\n  [\"73\u00d793=6789\", \"84\u00d793=7812\"],\n  [\"30\u00d756=1680\", \"68\u00d753=3604\"],\n  [\"66\u00d712=792\", \"92\u00d789=8188\"],\n  [\"87\u00d791=7917\", \"31\u00d728=868\"],\n  [\"94\u00d789=8366\", \"17\u00d798=1666\"],\n  [\"55\u00d733=1815\", \"32\u00d774=2368\"],\n  [\"33\u00d727=891\", \"89\u00d791=8099\"],\n  [\"51\u00d787=4437\", \"55\u00d753=2915\"],\n  [\"29\u00d740=1160\", \"25\u00d792=2300\"],\n  [\"71\u00d785=6035\", \"12\u00d722=264\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication problem text in the\n# document body with its updated value, as described by the diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"20\u00d765=1300\", \"21\u00d776=1596\"),\n  @(\"99\u00d785=8415\", \"54\u00d791=4914\"),\n  @(\"59\u00d790=5310\", \"38\u00d787=3306\"),\n  @(\"15\u00d751=765\",  \"61\u00d758=3538\"),\n  @(\"89\u00d787=7743\", \"43\u00d714=602\"),\n  @(\"96\u00d758=5568\", \"36\u00d748=1728\"),\n  @(\"22\u00d791=2002\", \"25\u00d798=2450\"),\n  @(\"42\u00d797=4074\", \"43\u00d785=3655\"),\n  @(\"25\u00d749=1225\", \"24\u00d722=528\"),\n  @(\"27\u00d774=1998\", \"25\u00d716=400\"),\n  @(\"29\u00d759=1711\", \"74\u00d790=6660\"),\n  @(\"94\u00d719=1786\", \"27\u00d731=837\"),\n  @(\"21\u00d797=2037\", \"13\u00d753=689\"),\n  @(\"65\u00d735=2275\", \"91\u00d731=2821\"),\n  @(\"63\u00d795=5985\", \"65\u00d775=4875\"),\n  @(\"73\u00d793=6789\", \"84\u00d793=7812\"),\n  @(\"30\u00d756=1680\", \"68\u00d753=3604\"),\n  @(\"66\u00d712=792\",  \"92\u00d789=8188\"),\n  @(\"87\u00d791=7917\", \"31\u00d728=868\"),\n  @(\"94\u00d789=8366\", \"17\u00d798=1666\"),\n  @(\"55\u00d733=1815\", \"32\u00d774=2368\"),\n  @(\"33\u00d727=891\",  \"89\u00d791=8099\"),\n  @(\"51\u00d787=4437\", \"55\u00d753=2915\"),\n  @(\"29\u00d740=1160\", \"25\u00d792=2300\"),\n  @(\"71\u00d785=6035\", \"12\u00d722=264\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n"}
